$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOPICS")
$ws.Activate()

# Add "Serverless Fns" in A36, matching the formatting used by the
# preceding entries (row 34's style: centered, wrap text, grey fill).
$ws.Range("A34").Copy()
$ws.Range("A36").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A36").Value = "Serverless Fns"

# Add "Hono" in A38, same formatting.
$ws.Range("A34").Copy()
$ws.Range("A38").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A38").Value = "Hono"

$excel.CutCopyMode = $false

# Match the author's final scroll position / selection in the sheet view.
$excel.ActiveWindow.ScrollRow = 31
$ws.Range("H37").Select()
